$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Archetype Name","Class","Identifier 1","Hash Code 1","Identifier 2","Hash Code 2","Anti Identifier","Hash Code 3"),
    @("Anomaly Forest","Forest","Wellspring Elf Princess","6lXRo","Life Banquet","6-rAw","None","None"),
    @("Amataz Forest","Forest","Amataz, Fairy Blader","6pQTI","Divine Smithing","6pOlw","None","None"),
    @("Pure Terror Forest","Forest","Terrorformer","6-suS","Colloseum on High","6sp_A","None","None"),
    @("Roach Accel Forest","Forest","Resolve of the Fallen","6-SmQ","Whirlwind Rhinoceroach","6lZu2","Fertile Aether","6pMJg"),
    @("Roach Natura Forest","Forest","Whirlwind Rhinoceroach","6lZu2","Fertile Aether","6pMJg","None","None"),
    @("Evo Forest","Forest","Zeus, The Supreme","6lDvy","Phantasmal Fairy Dragon","6tEcc","None","None"),
    @("Midrange Sword","Sword","Kagemitsu, Matchless Blade","6td0o","Kagemitsu, Matchless Blade","6td0o","None","None"),
    @("Dirt Rune","Rune","Forbidden Darkmage","6t_Rc","Karyl, Catty Sorceress","6t_RS","Riley, Hydroshaman","6q8s2"),
    @("Item Shop Rune","Rune","Arcane Item Shop","6q95g","Arcane Item Shop","6q95g","None","None"),
    @("Spellboost Rune","Rune","Kuon, Founder of Onmyodo","6t_RI","Runie, Resolute Diviner","6_djc","None","None"),
    @("Natura Rune","Rune","Riley, Hydroshaman","6q8s2","Karyl, Catty Sorceress","6t_RS","Magical Augmentation","6_X7Q"),
    @("Whale Dragon","Dragon","Eternal Whale","6_zhi","Slaughtering Dragonewt","6_zhY","None","None"),
    @("Discard Dragon","Dragon","Shipsbane Plesiosaurus","6yB-y","Turncoat Dragon Summoner","7007y","None","None"),
    @("Natura Dragon","Dragon","Valdain, Cursed Shadow","6qZiy","Whirlwind Pteranodon","6qUqS","None","None"),
    @("Jerva Evo Dragon","Dragon","Jerva, Wyrm Transcendent","6yB_6","Jerva, Wyrm Transcendent","6yB_6","None","None"),
    @("Showdown Dragon","Dragon","Titanic Showdown","6-Q4Y","Erntz, Justice","7007o","None","None"),
    @("Nepthys Shadow","Shadow","Nephthys, Goddess of Amenta","70OYc","Nephthys, Goddess of Amenta","70OYc","None","None"),
    @("Burial Rite Shadow","Shadow","He Who Once Rocked","6yaPS","Fatal Order","gXuCg","None","None"),
    @("Burial Rite Shadow","Shadow","He Who Once Rocked","6yaPS","Fatal Order","6yWFg","None","None"),
    @("PTP Shadow","Shadow","Hades, Father of Purgatory","6n7-I","Gremory, Death Teller","6yaPI","None","None"),
    @("Midrange Shadow","Shadow","Shuten-Doji","6ujqC","Ginsetsu, Great Fox","6umGI","None","None"),
    @("Natura Shadow","Shadow","Lubelle, Necrofamily","6qy7I","Thoth","6qy7S","None","None"),
    @("Control Blood","Blood","XII. Zelgenea, The World","6-UTo","Nerea, Beast Empress","6yypo","Permafrost Behemoth","None"),
    @("Corruption Blood","Blood","Baal","70mz6","Ravening Corruption","6ywNi","Contemptous Demon","None"),
    @("Wrath Blood","Blood","Permafrost Behemoth","6v8h6","Contemptous Demon","70kWi","None","None"),
    @("Satan Haven","Haven","Temple of Heresy","717QI","Temple of Heresy","717QI","None","None"),
    @("Elana Haven","Haven","Elana, Purest Prayer","6nupS","Elana, Purest Prayer","6nupS","None","None"),
    @("Sofina Amulet Haven","Haven","VIII. Sofina, Strength","719NI","Master Adjudicator","6rgW2","None","None"),
    @("Ward Haven","Haven","Wilbert, Grand Knight","6vX5I","Wilbert, Grand Knight","6vX5I","None","None"),
    @("Ra Control Haven","Haven","Ra, Radiance Incarnate","719Nc","Ra, Radiance Incarnate","719Nc","None","None"),
    @("Artifact Portal","Portal","Artifact Scan","71TeA","Absolute Modesty","6s5My","None","None"),
    @("Artifact Portal","Portal","Artifact Scan","gZ1SK","Absolute Modesty","6s5My","None","None"),
    @("Belphomet Portal","Portal","Worldreaver, Belphomet","6zjey","Aerial Craft","6oEni","None","None"),
    @("Orchis Portal","Portal","Orchis, Linked Heart","6oHE6","Orchis, Linked Heart","6oHE6","None","None"),
    @("Control Portal","Portal","Karula, Arts Master","6vvVo","XII. Zelgenea, The World","6-UTo","None","None")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i+1, $j+1).Value = $row[$j]
    }
}

$ws.Range("H26").Select()

Write-Host "Updated fortune hand meta sheet"